$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$valuesA = @(
    -0.221225574440636,
    -0.4529329355182927,
    -0.5289049001564117,
    -0.1867077086173133,
    -0.1425258833026161,
    -0.2486432360548048,
    -0.4821868847663077,
    -0.6444530289452352,
    -0.4216969412043984,
    -0.2451629587420543,
    -0.1677887243925625,
    -0.1024247910263559,
    -0.2323144740892847,
    -0.07433902824470077,
    -0.28644912950437,
    -0.06798365622505377,
    0.01490435119564006,
    0.03643196019040799,
    -0.1778476170974412,
    -0.008985210202750815,
    -0.07236841886345335,
    0.06300908772083888,
    0.9000937714275231,
    0.1350499628996302,
    0.1547467607237886,
    0.1425981655686474,
    0.5019926804822898,
    0.7394975758911319,
    0.2216424931757903,
    0.1115071214070912,
    0.2284796900786265,
    0.1272358686940622,
    0.1075420983023397,
    0.3758252025780184,
    0.2226193700633499,
    -0.03122971852243707,
    0.3454078674352573,
    -0.08797232348903913,
    0.09899929987938683,
    -0.09493465746259849,
    0.2664362332282849,
    0.04451977215206121,
    0.007223945691946815,
    -0.1433166457959914,
    -0.1662375939550853,
    -0.1647295940943758,
    -0.2179105188753415,
    -0.226217399964751,
    -0.09145167631939459,
    -0.2365181307787485,
    -0.2365181307787485,
    -0.1899107247446276,
    -0.2177731319697336,
    -0.1532242862750103,
    -0.1105641561091033,
    -0.1820161885290238,
    -0.1244150338271487,
    -0.2070706678073113,
    -0.2521507862176968,
    -0.2628710625733858,
    -0.1492751247279001,
    -0.4817226268554734,
    -0.2850620294950607,
    -0.2456735080486877,
    -0.08101856810283901,
    0.02665950280576842
)

$valuesB = @(
    -0.2987113970417993,
    -0.3882324487736706,
    -0.468725384260669,
    -0.2176596936741573,
    -0.1165560554551552,
    -0.2813035583867065,
    -0.4097452871424743,
    -0.5068947078914234,
    -0.2167961031478336,
    -0.2682052297001176,
    -0.1592762765478757,
    -0.09228825415552418,
    -0.1642208691978131,
    -0.01562184063545382,
    -0.1320365665804203,
    0.06552326622987384,
    0.04479460899560108,
    0.08336493776933186,
    -0.1059617328113978,
    0.06165242795151016,
    0.07275244515004033,
    0.07093169640599929,
    0.7371844791573458,
    0.107947683635579,
    0.1387981892610986,
    0.142316270878525,
    0.332523306805975,
    0.6120116303542912,
    0.1815843056231921,
    0.08800195546024515,
    0.2157938649752954,
    0.1329301324753304,
    0.0773005947336623,
    0.2752186466059716,
    0.1034232283275243,
    -0.03240632934089572,
    0.2855184078394947,
    -0.139589759800075,
    0.1370135359703478,
    -0.1673046903933172,
    0.2667067070515314,
    0.08201273402688293,
    -0.05641184709164448,
    -0.1034158618095164,
    -0.1543678915519102,
    -0.151560487986453,
    -0.1926865224680253,
    -0.1948180542766703,
    -0.1006449971721333,
    -0.2181183661709926,
    -0.2181183661709926,
    -0.1719949621505628,
    -0.1883582628345376,
    -0.1339516217372279,
    -0.1029230541526095,
    -0.1396509536707337,
    -0.1657163046327254,
    -0.2158717456268029,
    -0.2434085372046686,
    -0.2138441324529243,
    -0.09111869749647811,
    -0.4405003228457873,
    -0.2660528363680922,
    -0.2288769295655513,
    -0.06736570766589559,
    0.005499169068905142
)

for ($i = 0; $i -lt $valuesA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $valuesA[$i]
    $ws.Cells.Item($row, 2).Value = $valuesB[$i]
}